$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(5).ColumnWidth = 28.5703125
$ws.Columns.Item(5).HorizontalAlignment = -4131
$ws.Columns.Item(5).VerticalAlignment = -4108
$ws.Columns.Item(5).IndentLevel = 1

$ws.Range("E1").Font.Bold = $true
